# Add two new columns, I ("I0") and J ("IF"), to the right of the
# existing data table (which currently ends at column H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Cells.Item(1, 9).Value  = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Match the header formatting used by the rest of row 1 (bold, centered,
# bordered) by copying H1's format onto the two new header cells.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (2-15) -------------------------------------------------
$values = @{
  2  = @(8, 9)
  3  = @(4, 5)
  4  = @(5, 7)
  5  = @(9, 9)
  6  = @(7, 7)
  7  = @(9, 9)
  8  = @(5, 6)
  9  = @(7, 7)
  10 = @(9, 9)
  11 = @(7, 7)
  12 = @(7, 7)
  13 = @(8, 8)
  14 = @(4, 4)
  15 = @(6, 6)
}

foreach ($row in 2..15) {
  $pair = $values[$row]
  $ws.Cells.Item($row, 9).Value  = $pair[0]
  $ws.Cells.Item($row, 10).Value = $pair[1]
}
